$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.237.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.72%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.656.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.64%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.75%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5226"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.47%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.73%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06332"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.41%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.96%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.553"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.23%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.648.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.884.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.02%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5675"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8123"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.16%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'65.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.63%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.232.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.64%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'192.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.040"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.85%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.63%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'143.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.277"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.84%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.497"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05607"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.280"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.77%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.45%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.380"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.586"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.41%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.66%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.9449"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.39%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5750"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.74%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01599"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.08%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.914"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.36%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.581"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8482"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.25%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.67%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.032.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.62%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'102.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.795.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'58.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -2.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05317"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4355"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.20%  "
$ws.Range("E51").Style = "Normal"
